# Commit: "Fruta / hortaliza, semanal"
# A new weekly record is inserted as row 130, pushing the existing rows
# 130-236 down to 131-237 (data unchanged, just shifted down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 130; this shifts rows 130-236 down to 131-237
# and keeps all their existing values/styles intact.
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with the new record's data.
$ws.Range("A130").Value = 3
$ws.Range("B130").Value = "Femacal de La Calera"
$ws.Range("C130").Value = "Coquimbo"
$ws.Range("D130").Value = 44977
$ws.Range("E130").Value = 5
$ws.Range("F130").Value = 100112052
$ws.Range("G130").Value = "Albahaca"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 105
$ws.Range("K130").Value = 5000
$ws.Range("L130").Value = 5500
$ws.Range("M130").Value = 5238
$ws.Range("N130").Value = "`$/docena de matas"
$ws.Range("O130").Value = "Provincia de Quillota"
$ws.Range("P130").Value = 873
$ws.Range("Q130").Value = 6
$ws.Range("R130").Value = "Hortaliza"
